$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ViewProductThroughList")

$ws.Rows("3").Copy() | Out-Null
$ws.Rows("4").Insert() | Out-Null

$ws.Range("D4").Value = "Máy nước nóng Kangaroo 22 lít KG 70A2"
$ws.Range("E4").Value = "PASSED"
$ws.Range("F4").Value = "Máy nước nóng"
$ws.Range("G4").Value = "Kangaroo"
$ws.Range("H4").Value = "Máy nước nóng Kangaroo 2500W KG 70A2"

$ws.Range("E4").Select()
